$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-09-09 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-10 Sunday", 2) | Out-Null

# Update each math-problem cell in the table by position (row, col)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "33+57="
$t.Cell(1, 2).Range.Text = "80+13="
$t.Cell(1, 3).Range.Text = "53+14="
$t.Cell(1, 4).Range.Text = "67+12="
$t.Cell(1, 5).Range.Text = "84-14="
$t.Cell(2, 1).Range.Text = "75-43="
$t.Cell(2, 2).Range.Text = "31+19="
$t.Cell(2, 3).Range.Text = "68-43="
$t.Cell(2, 4).Range.Text = "97-77="
$t.Cell(2, 5).Range.Text = "79-79="
$t.Cell(3, 1).Range.Text = "36-13="
$t.Cell(3, 2).Range.Text = "9+48="
$t.Cell(3, 3).Range.Text = "47+34="
$t.Cell(3, 4).Range.Text = "31-7="
$t.Cell(3, 5).Range.Text = "20+41="
$t.Cell(4, 1).Range.Text = "61-11="
$t.Cell(4, 2).Range.Text = "60-14="
$t.Cell(4, 3).Range.Text = "54-11="
$t.Cell(4, 4).Range.Text = "74-23="
$t.Cell(4, 5).Range.Text = "73-55="
$t.Cell(5, 1).Range.Text = "58-31="
$t.Cell(5, 2).Range.Text = "68-67="
$t.Cell(5, 3).Range.Text = "53+10="
$t.Cell(5, 4).Range.Text = "44-5="
$t.Cell(5, 5).Range.Text = "64-3="
$t.Cell(6, 1).Range.Text = "66-16="
$t.Cell(6, 2).Range.Text = "11+29="
$t.Cell(6, 3).Range.Text = "6+15="
$t.Cell(6, 4).Range.Text = "55+42="
$t.Cell(6, 5).Range.Text = "27+1="
$t.Cell(7, 1).Range.Text = "67-11="
$t.Cell(7, 2).Range.Text = "23-2="
$t.Cell(7, 3).Range.Text = "94-70="
$t.Cell(7, 4).Range.Text = "70-13="
$t.Cell(7, 5).Range.Text = "11+42="
$t.Cell(8, 1).Range.Text = "8+20="
$t.Cell(8, 2).Range.Text = "14+23="
$t.Cell(8, 3).Range.Text = "21+65="
$t.Cell(8, 4).Range.Text = "72-56="
$t.Cell(8, 5).Range.Text = "77-13="
$t.Cell(9, 1).Range.Text = "25+19="
$t.Cell(9, 2).Range.Text = "86-79="
$t.Cell(9, 3).Range.Text = "52+36="
$t.Cell(9, 4).Range.Text = "15+43="
$t.Cell(9, 5).Range.Text = "38-2="
$t.Cell(10, 1).Range.Text = "21+51="
$t.Cell(10, 2).Range.Text = "23+75="
$t.Cell(10, 3).Range.Text = "76+0="
$t.Cell(10, 4).Range.Text = "54-20="
$t.Cell(10, 5).Range.Text = "73-1="
$t.Cell(11, 1).Range.Text = "33+46="
$t.Cell(11, 2).Range.Text = "15+38="
$t.Cell(11, 3).Range.Text = "79-25="
$t.Cell(11, 4).Range.Text = "40+55="
$t.Cell(11, 5).Range.Text = "6+87="
$t.Cell(12, 1).Range.Text = "18-14="
$t.Cell(12, 2).Range.Text = "75+10="
$t.Cell(12, 3).Range.Text = "27+3="
$t.Cell(12, 4).Range.Text = "35-5="
$t.Cell(12, 5).Range.Text = "70-50="
$t.Cell(13, 1).Range.Text = "99-56="
$t.Cell(13, 2).Range.Text = "70-13="
$t.Cell(13, 3).Range.Text = "24+69="
$t.Cell(13, 4).Range.Text = "92-31="
$t.Cell(13, 5).Range.Text = "47+2="
$t.Cell(14, 1).Range.Text = "72-38="
$t.Cell(14, 2).Range.Text = "91-39="
$t.Cell(14, 3).Range.Text = "7+50="
$t.Cell(14, 4).Range.Text = "24+53="
$t.Cell(14, 5).Range.Text = "82-31="
$t.Cell(15, 1).Range.Text = "40+31="
$t.Cell(15, 2).Range.Text = "44+38="
$t.Cell(15, 3).Range.Text = "79+1="
$t.Cell(15, 4).Range.Text = "13+38="
$t.Cell(15, 5).Range.Text = "39+14="
$t.Cell(16, 1).Range.Text = "53-6="
$t.Cell(16, 2).Range.Text = "66+17="
$t.Cell(16, 3).Range.Text = "33-15="
$t.Cell(16, 4).Range.Text = "74+2="
$t.Cell(16, 5).Range.Text = "5+9="
$t.Cell(17, 1).Range.Text = "34+48="
$t.Cell(17, 2).Range.Text = "86-62="
$t.Cell(17, 3).Range.Text = "85+5="
$t.Cell(17, 4).Range.Text = "71-43="
$t.Cell(17, 5).Range.Text = "71-59="
$t.Cell(18, 1).Range.Text = "79-47="
$t.Cell(18, 2).Range.Text = "61-28="
$t.Cell(18, 3).Range.Text = "40-33="
$t.Cell(18, 4).Range.Text = "37+43="
$t.Cell(18, 5).Range.Text = "39+19="
$t.Cell(19, 1).Range.Text = "29+45="
$t.Cell(19, 2).Range.Text = "17+14="
$t.Cell(19, 3).Range.Text = "53+4="
$t.Cell(19, 4).Range.Text = "3+18="
$t.Cell(19, 5).Range.Text = "0+40="
$t.Cell(20, 1).Range.Text = "77-37="
$t.Cell(20, 2).Range.Text = "32+19="
$t.Cell(20, 3).Range.Text = "29+49="
$t.Cell(20, 4).Range.Text = "95-84="
$t.Cell(20, 5).Range.Text = "24+18="
